# This script reproduces the commit that removes the AQL "id" parameter
# (x="'value1'") from the ":userdoc" field instruction, and splits the
# "User document part Texte 1" run in two (inserting a _GoBack bookmark
# at the edit point), matching the author's manual edits in Word.
#
# Field instruction text (instrText) and bookmarks cannot be reliably
# rewritten through Range/Selection text assignment in this host (the
# field-code span collapses to a single anchor point), so instead we
# replace the whole target paragraph's XML using Range.InsertXML with a
# minimal WordProcessingML package, which lets us express the exact
# target run/bookmark structure.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document ' + $wNs + '><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData>' +
           '</pkg:part></pkg:package>'
    $r.InsertXML($pkg) | Out-Null
}

# --- Paragraph 2: the " m:userdoc x='value1' " field -> " m:userdoc    " ---
$fieldPara = $d.Paragraphs.Item(2)
$fieldCode = $d.Fields.Item(1).Code.Text
if ($fieldCode -notmatch "userdoc") {
    throw "Unexpected field code, aborting: $fieldCode"
}

$p2xml = '<w:p w:rsidR="00A7781B" w:rsidRDefault="005C73CF" w:rsidP="005C73CF">' +
         '<w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr>' +
         '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
         '<w:r><w:instrText xml:space="preserve"> </w:instrText></w:r>' +
         '<w:r w:rsidR="0056766F"><w:instrText>m</w:instrText></w:r>' +
         '<w:r w:rsidR="00A7781B"><w:instrText>:userdoc</w:instrText></w:r>' +
         '<w:r><w:instrText xml:space="preserve"> </w:instrText></w:r>' +
         '<w:r><w:instrText xml:space="preserve">   </w:instrText></w:r>' +
         '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
         '</w:p>'

Replace-ParagraphXml 2 $p2xml

# --- Paragraph 3: "User document part Texte 1" -> "User document" + _GoBack + " part Texte 1" ---
$textPara = $d.Paragraphs.Item(3)
if ($textPara.Range.Text -notmatch "User document part Texte 1") {
    throw "Unexpected paragraph text, aborting: $($textPara.Range.Text)"
}

$p3xml = '<w:p w:rsidR="00A7781B" w:rsidRDefault="00A7781B" w:rsidP="005C73CF">' +
         '<w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr>' +
         '<w:r><w:t>User document</w:t></w:r>' +
         '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
         '<w:bookmarkEnd w:id="0"/>' +
         '<w:r><w:t xml:space="preserve"> part Texte 1</w:t></w:r>' +
         '</w:p>'

Replace-ParagraphXml 3 $p3xml
